# Renaming placeholder stream names
# Updates unit-sheet tab names and the "placeholder" stream-name cells
# (B1 titles, inlet/outlet "stream" columns, and two species-list cells)
# so they use the new F53.. numbering instead of the old F101/F102/etc.
# placeholders.

$wb = $excel.ActiveWorkbook

# --- Cell-content edits, keyed by the sheet's ORIGINAL (pre-rename) name ---
# Each entry: sheet name -> list of (cell address, new value)
$cellEdits = @{
    "X101_TBP_Extraction" = @(
        @{ Cell = "B19"; Value = "Sm(NO3)3, Sr(NO3)2, CsNO3, HNO3, I_aq, Gd(NO3)3, H2O, Eu(NO3)3, Nd(NO3)3" }
    )
    "CF101_Coalescer_X101_to_E101" = @(
        @{ Cell = "B34"; Value = "F53" }
        @{ Cell = "B35"; Value = "F54" }
    )
    "HX__F101-T" = @(
        @{ Cell = "B1";  Value = "HX__F53-T" }
        @{ Cell = "B15"; Value = "F53" }
        @{ Cell = "B19"; Value = "F53-T" }
    )
    "PC__F101-p" = @(
        @{ Cell = "B1";  Value = "PC__F53-p" }
        @{ Cell = "B15"; Value = "F53-T" }
        @{ Cell = "B19"; Value = "F53-p" }
    )
    "E101_Evaporator" = @(
        @{ Cell = "B38"; Value = "F53-p" }
    )
    "X102_AHA_Strip" = @(
        @{ Cell = "B18"; Value = "AHA, H2O, HNO3" }
    )
    "CF102_Coalescer_X102_to_E102" = @(
        @{ Cell = "B34"; Value = "F55" }
        @{ Cell = "B35"; Value = "F56" }
    )
    "HX__F102-T" = @(
        @{ Cell = "B1";  Value = "HX__F55-T" }
        @{ Cell = "B15"; Value = "F55" }
        @{ Cell = "B19"; Value = "F55-T" }
    )
    "PC__F102-p" = @(
        @{ Cell = "B1";  Value = "PC__F55-p" }
        @{ Cell = "B15"; Value = "F55-T" }
        @{ Cell = "B19"; Value = "F55-p" }
    )
    "E102_Evaporator" = @(
        @{ Cell = "B38"; Value = "F55-p" }
    )
    "CF103_Coalescer_X103_to_E103" = @(
        @{ Cell = "B34"; Value = "F57" }
        @{ Cell = "B35"; Value = "F58" }
    )
    "HX__F112-T" = @(
        @{ Cell = "B1";  Value = "HX__F57-T" }
        @{ Cell = "B15"; Value = "F57" }
        @{ Cell = "B19"; Value = "F57-T" }
    )
    "E103_Evaporator100C" = @(
        @{ Cell = "B19"; Value = "F57-T" }
    )
    "V104_AcidPurge" = @(
        @{ Cell = "B19"; Value = "F59" }
    )
    "V105_KOWaterPurge" = @(
        @{ Cell = "B19"; Value = "F60" }
    )
    "M132_CoalescerOrgMixer" = @(
        @{ Cell = "B14"; Value = "F54" }
        @{ Cell = "B15"; Value = "F56" }
        @{ Cell = "B16"; Value = "F58" }
        @{ Cell = "B20"; Value = "F61" }
    )
    "M133_OrgRecoveryMixer" = @(
        @{ Cell = "B15"; Value = "F61" }
        @{ Cell = "B19"; Value = "F62" }
    )
    "HX133_SolventCooldown" = @(
        @{ Cell = "B15"; Value = "F62" }
    )
}

foreach ($sheetName in $cellEdits.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($edit in $cellEdits[$sheetName]) {
        $ws.Range($edit.Cell).Value = $edit.Value
    }
}

# --- Worksheet tab renames (done after the cell edits so the lookups
#     above can still use the original tab names) ---
$tabRenames = @(
    @{ Old = "HX__F101-T"; New = "HX__F53-T" }
    @{ Old = "PC__F101-p"; New = "PC__F53-p" }
    @{ Old = "HX__F102-T"; New = "HX__F55-T" }
    @{ Old = "PC__F102-p"; New = "PC__F55-p" }
    @{ Old = "HX__F112-T"; New = "HX__F57-T" }
)

foreach ($rename in $tabRenames) {
    $ws = $wb.Worksheets.Item($rename.Old)
    $ws.Name = $rename.New
}
